$d = $word.ActiveDocument

$replacements = @(
    @("397÷8=49, 5", "168÷9=18, 6"),
    @("492÷5=98, 2", "313÷9=34, 7"),
    @("299÷3=99, 2", "343÷7=49, 0"),
    @("487÷3=162, 1", "385÷3=128, 1"),
    @("978÷9=108, 6", "262÷8=32, 6"),
    @("742÷7=106, 0", "959÷6=159, 5"),
    @("846÷9=94, 0", "804÷2=402, 0"),
    @("909÷8=113, 5", "756÷7=108, 0"),
    @("897÷3=299, 0", "790÷4=197, 2"),
    @("669÷6=111, 3", "355÷8=44, 3"),
    @("213÷7=30, 3", "431÷3=143, 2"),
    @("628÷5=125, 3", "364÷2=182, 0"),
    @("413÷7=59, 0", "425÷6=70, 5"),
    @("261÷8=32, 5", "827÷8=103, 3"),
    @("485÷8=60, 5", "215÷2=107, 1"),
    @("927÷4=231, 3", "792÷8=99, 0"),
    @("992÷3=330, 2", "187÷3=62, 1"),
    @("750÷6=125, 0", "241÷3=80, 1"),
    @("693÷3=231, 0", "645÷6=107, 3"),
    @("633÷8=79, 1", "929÷7=132, 5"),
    @("611÷8=76, 3", "833÷8=104, 1"),
    @("460÷6=76, 4", "288÷9=32, 0"),
    @("151÷7=21, 4", "797÷5=159, 2"),
    @("570÷3=190, 0", "158÷8=19, 6"),
    @("191÷2=95, 1", "588÷3=196, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
